# Remove Wind (onshore) / Wind (offshore) related rows & columns from the
# urbs classic input workbook, matching the commit
# "deleted offshore and onshore from urbs classic inputs".

$wb = $excel.ActiveWorkbook

# --- Commodity sheet: delete rows for WindOff / WindOn commodities -------
$wsCommodity = $wb.Worksheets.Item("Commodity")
$wsCommodity.Rows("2:3").Delete()

# --- Process sheet: delete rows for Wind (onshore) / Wind (offshore) -----
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcess.Rows("2:3").Delete()

# --- Process-Commodity sheet: delete the 4 rows describing the In/Out ----
# relationships of Wind (onshore) and Wind (offshore)
$wsProcessCommodity = $wb.Worksheets.Item("Process-Commodity")
$wsProcessCommodity.Rows("2:5").Delete()

# --- SupIm sheet: delete the WindOn / WindOff capacity factor columns ----
$wsSupIm = $wb.Worksheets.Item("SupIm")
$wsSupIm.Columns("B:C").Delete()

# Make SupIm the active sheet/selection, matching the saved view state.
$wsSupIm.Activate()
$wsSupIm.Columns("B:C").Select()
